$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update separators (comma -> semicolon) and synonym lists in DictionaryColNames (B) and SIMSmethods (C) columns ---
$ws.Range("B2").Value = 'File; Filename'
$ws.Range("B3").Value = 'Comment; Comment..2'
$ws.Range("B4").Value = '\u03B418O \u2030 VSMOW vs UWC-3; d18O \u2030 VSMOW; d18_VSMOW; d18O [VSMOW],d18 VSMOW; \u03B418O \u2030 VSMOW; <U+03B4>18O <U+2030> VSMOW vs UWC-3; δ18O ‰ VSMOW vs UWC-3; d18O ‰ VSMOW; δ18O ‰ VSMOW'
$ws.Range("B5").Value = '2SD (ext.); Er (2S); 2SD; Std_1SD; Er(2S)'
$ws.Range("B6").Value = 'Bias; IMF; Mass Bias (<U+2030>); Mass Bias (‰)'
$ws.Range("B7").Value = 'd18O ‰ raw; d18O_m; d18O_meas; d18_c; d18O meas; d18O ‰ measured; \u03B418O ‰ measured; <U+03B4>18O <U+2030> measured; δ18O ‰ measured'
$ws.Range("B8").Value = '2SE (int.); d18O-2SE; d13C-2SE'
$ws.Range("B9").Value = '16O (Gcps); 16O(E9 cps); 16O     (E9 cps); 16O (E9); 16O(E9)'
$ws.Range("B10").Value = 'IP(nA); IP(nA)  1.7 to 1.9; IP (nA)'
$ws.Range("B12").Value = 'Yield (Gcps/nA); Yield(E9cps/nA); Yield (E9cps/nA)'
$ws.Range("B13").Value = 'date; Date'
$ws.Range("B14").Value = 'time; Time'
$ws.Range("B15").Value = 'X; x'
$ws.Range("B16").Value = 'Y; y'
$ws.Range("B17").Value = 'DTFA-X; DTFA X'
$ws.Range("B18").Value = 'DTFA-Y; DTFA Y'
$ws.Range("B19").Value = '16OH/16O; 16O1H/16O; 13CH/13C'
$ws.Range("B20").Value = '12C (E6); 12C(E6)'
$ws.Range("B21").Value = 'd13_C; δ13C ‰ VPDB; d13C PDB; δ13C [‰, PDB]; \u03B413C [‰ PDB]; δ13C [‰, VPDB]; \u03B413C [‰, VPDB]; \u03B413C [\u2030, VPDB]; \u03B413C [\u2030, PDB]; ; \u03B413C [\u2030 PDB]'
$ws.Range("B22").Value = 'd13C_m; δ13C ‰ measured; \u03B413C measured; \u03B413C \u2030 measured'

$ws.Range("C2").Value = 'd18O10; d13C7'
$ws.Range("C3").Value = 'd18O10; d13C7'
$ws.Range("C5").Value = 'd18O10; d13C7'
$ws.Range("C8").Value = 'd18O10; d13C7'
$ws.Range("C10").Value = 'd18O10; d13C7'
$ws.Range("C13").Value = 'd18O10; d13C7'
$ws.Range("C14").Value = 'd18O10; d13C7'
$ws.Range("C15").Value = 'd18O10; d13C7'
$ws.Range("C16").Value = 'd18O10; d13C7'
$ws.Range("C17").Value = 'd18O10; d13C7'
$ws.Range("C18").Value = 'd18O10; d13C7'
$ws.Range("C19").Value = 'd18O10; d13C7'

# --- Insert new row for "Mass" column dictionary entry (new row 19, pushes Hyd..INDEX down by one) ---
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = 'Mass'
$ws.Range("B19").Value = 'Mass; mass'
$ws.Range("C19").Value = 'd18O10; d13C7'
$ws.Range("D19").Value = 'bits'
$ws.Range("E19").Value = 'Numeric'

# --- Update selection to match final state ---
$ws.Range("B23").Select()

# --- Update window size (workbookView) ---
$excel.ActiveWindow.Width = 24360
$excel.ActiveWindow.Height = 15480
